# Testdata.xlsx update: "updated the pagebase removed getdriver method"
#
# 1. Rename sheet "Locators" -> "HomePageLocators"
# 2. Add a new column G to Prod_TD with a header/value pair copied from the
#    HomePageLocators sheet (HomeSearchLocator / div.search-icon)
# 3. Update the remembered selection on both sheets

$wb = $excel.ActiveWorkbook

$prod = $wb.Worksheets.Item("Prod_TD")
$loc  = $wb.Worksheets.Item("Locators")

# --- 1. Rename the Locators sheet ---------------------------------------
$loc.Name = "HomePageLocators"

# --- 2. Add column G (header + value), matching HomePageLocators A1/A2 --
# Copy formatting (style) from the existing header cell A1 so the new
# header cell G1 picks up the same style (font/no-fill) used by the rest
# of row 1.
$prod.Range("A1").Copy() | Out-Null
$prod.Range("G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$prod.Range("G1").Value = $loc.Range("A1").Value2   # "HomeSearchLocator"
$prod.Range("G2").Value = $loc.Range("A2").Value2   # "div.search-icon"

# Match the column width used on HomePageLocators column A (~18.86 chars).
$prod.Columns.Item(7).ColumnWidth = 18

$excel.CutCopyMode = $false

# --- 3. Update selections -------------------------------------------------
# HomePageLocators: select A1:A2, then hand focus back to Prod_TD so the
# active tab stays Prod_TD.
$loc.Range("A1:A2").Select() | Out-Null
$prod.Activate() | Out-Null
$prod.Range("D15").Select() | Out-Null

Write-Host ("Sheets: " + (($wb.Worksheets | ForEach-Object { $_.Name }) -join ", "))
